# benchmark cigre residential/commercial/industrial added
#
# Replaces the hourly weekday/weekend/holiday load-profile coefficients on
# sheet "1" (commercial CIGRE profile) with the new benchmark figures, and
# updates the saved UI selection state (active cell / active sheet) on a
# few sheets to match the author's final session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "1": new per-hour coefficients (same value repeated across the
# Monday..Holiday columns B:I) replacing the old placeholder 1 / 1.2 / 0.1
# figures, and the old centred number style (s="3") is cleared so the
# cells fall back to the default/general style - matching the target
# workbook exactly.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1")
$ws1.Activate()

$newValues = @{
    2  = 20.99
    3  = 20.25
    4  = 20.87
    5  = 21.74
    6  = 23.47
    7  = 25.46
    8  = 36.12
    9  = 49.27
    10 = 64.9
    11 = 79.41
    12 = 84.99
    13 = 89.83
    14 = 92.31
    15 = 91.94
    16 = 87.6
    17 = 85.61
    18 = 89.83
    19 = 99.75
    20 = 79.66
    21 = 69.98
    22 = 59.19
    23 = 48.9
    24 = 29.92
    25 = 19.75
}

foreach ($row in 2..25) {
    $rng = $ws1.Range("B" + $row + ":I" + $row)
    $rng.Value = $newValues[$row]
    $rng.ClearFormats()
}

$ws1.Range("J36").Select()

# ---------------------------------------------------------------------
# Sheet "2": only the remembered selection changes (same B2:I25 block).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2")
$ws2.Activate()
$ws2.Range("B2:I25").Select()

# ---------------------------------------------------------------------
# Sheet "6": loses the "last active" flag and its remembered selection
# becomes the single cell F35.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("6")
$ws6.Activate()
$ws6.Range("F35").Select()

# ---------------------------------------------------------------------
# Sheet "12" ends up the active tab in the saved workbook (activeTab
# counts from 0, so the 12th/last sheet -> 11), keeping its K18 selection.
# Activating it last is what flips tabSelected/activeTab in the saved file.
# ---------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("12")
$ws12.Activate()
$ws12.Range("K18").Select()
